# Issue with Parallel execution of feature files
#
# - Rename "Sheet1" -> "LogIn"
# - Move the selection on the LogIn sheet from D17 -> A5

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Name = "LogIn"

$ws.Activate()
$ws.Range("A5").Select()
